$d = $word.ActiveDocument

$d.Content.Find.Execute("1013310111806", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1015707121800", 2)

$d.Content.Find.Execute(" 10 листопада 2018 р.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 7 грудня 2018 р.", 2)

$d.Content.Find.Execute("396965", $true, $false, $false, $false, $false,
                         $true, 1, $false, "009973", 2)

$d.Content.Find.Execute("δQn = мінус 50.0%; δQt = мінус 16.7%; δQmin = мінус 140.0%", $true, $false, $false, $false, $false,
                         $true, 1, $false, "δQt = 80.0%; δQmin = 65.0%", 2)
